# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with freshly scraped values, preserving each cell as literal text
# (matches the inline-string cells already used throughout the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force-store as text even when the string parses as a number
    # (e.g. "214.60"), then drop back to the default "Normal" style
    # so no stray number-format style sticks to the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '25.963.38'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.636.92'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  -0.49%  '
Set-TextValue $ws.Range("D5") '214.60'
$ws.Range("E5").Value = '  -0.41%  '
Set-TextValue $ws.Range("D6") '0.5089'
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  -0.39%  '
Set-TextValue $ws.Range("D8") '0.2564'
$ws.Range("E8").Value = '  -0.63%  '
Set-TextValue $ws.Range("D9") '0.06350'
$ws.Range("E9").Value = '  -0.77%  '
Set-TextValue $ws.Range("D10") '19.65'
$ws.Range("E10").Value = '  +0.23%  '
Set-TextValue $ws.Range("D11") '0.07748'
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").Value = '1.642.48'
$ws.Range("E13").Value = '  +0.50%  '
Set-TextValue $ws.Range("D14") '0.5433'
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("D15").Value = '0.0₅7708'
$ws.Range("E15").Value = '  -2.08%  '
Set-TextValue $ws.Range("D16") '63.97'
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("D17").Value = '25.988.82'
$ws.Range("E17").Value = '  +0.05%  '
$ws.Range("E18").Value = '  -0.39%  '
Set-TextValue $ws.Range("D19") '198.81'
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -0.74%  '
Set-TextValue $ws.Range("D22") '6.040'
$ws.Range("E22").Value = '  +0.70%  '
Set-TextValue $ws.Range("D23") '1.004'
$ws.Range("E23").Value = '  -0.23%  '
Set-TextValue $ws.Range("D24") '1.890'
$ws.Range("E24").Value = '  +1.23%  '
Set-TextValue $ws.Range("D25") '141.12'
$ws.Range("E25").Value = '  +0.28%  '
Set-TextValue $ws.Range("D26") '0.1201'
$ws.Range("E26").Value = '  +5.10%  '
Set-TextValue $ws.Range("D27") '6.824'
$ws.Range("E27").Value = '  -0.78%  '
Set-TextValue $ws.Range("D28") '15.57'
$ws.Range("E28").Value = '  -1.20%  '
Set-TextValue $ws.Range("D29") '1.232'
$ws.Range("E29").Value = '  -0.89%  '
Set-TextValue $ws.Range("D30") '0.04898'
$ws.Range("E30").Value = '  -2.61%  '
Set-TextValue $ws.Range("D31") '3.255'
$ws.Range("E31").Value = '  -0.38%  '
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("E33").Value = '  -0.45%  '
Set-TextValue $ws.Range("D34") '2.369'
$ws.Range("E34").Value = '  -0.08%  '
Set-TextValue $ws.Range("D35") '0.9075'
$ws.Range("E35").Value = '  +1.53%  '
Set-TextValue $ws.Range("D36") '2.584'
$ws.Range("E36").Value = '  -1.12%  '
$ws.Range("D37").Value = '1.126.48'
$ws.Range("E37").Value = '  -1.59%  '
Set-TextValue $ws.Range("D38") '0.5455'
$ws.Range("E38").Value = '  -1.71%  '
Set-TextValue $ws.Range("D39") '0.01560'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").Value = '  -0.45%  '
Set-TextValue $ws.Range("D41") '2.522'
$ws.Range("E41").Value = '  -1.59%  '
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("E43").Value = '  +2.30%  '
Set-TextValue $ws.Range("D44") '99.01'
$ws.Range("E44").Value = '  -0.86%  '
Set-TextValue $ws.Range("D45") '5.419'
$ws.Range("E45").Value = '  -4.83%  '
$ws.Range("D46").Value = '1.776.23'
$ws.Range("E46").Value = '  -0.31%  '
Set-TextValue $ws.Range("D47") '0.4524'
$ws.Range("E47").Value = '  -0.10%  '
Set-TextValue $ws.Range("D48") '1.006'
$ws.Range("E48").Value = '  -0.04%  '
Set-TextValue $ws.Range("D49") '54.90'
$ws.Range("E49").Value = '  -0.96%  '
Set-TextValue $ws.Range("D50") '0.05112'
$ws.Range("E50").Value = '  +1.02%  '
Set-TextValue $ws.Range("D51") '1.005'
$ws.Range("E51").Value = '  -0.10%  '
